$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = "Yes"
$ws.Range("C3").Value = "Yes"
$ws.Range("C4").Value = "Yes"
$ws.Range("C60").Value = "Yes"
$ws.Range("C61").Value = "Yes"
$ws.Range("C63").Value = "Yes"

$excel.ActiveWindow.ScrollRow = 34
$ws.Range("C4").Select() | Out-Null
